# NEL_deaths.xlsx update: refresh daily/cumulative death counts and append
# six new days of data (28 Mar 2020 - 2 Apr 2020 / serials 43945-43950).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Corrections to existing daily figures (column B = DailyBHRUT)
# ---------------------------------------------------------------------
$ws.Range("B6").Value  = 0
$ws.Range("B25").Value = 4
$ws.Range("B30").Value = 15
$ws.Range("B37").Value = 11
$ws.Range("B38").Value = 17
$ws.Range("B39").Value = 17
$ws.Range("B41").Value = 18
$ws.Range("B44").Value = 16
$ws.Range("B45").Value = 14
$ws.Range("B50").Value = 10
$ws.Range("B54").Value = 5
$ws.Range("F54").Value = 3
$ws.Range("B55").Value = 3
$ws.Range("F55").Value = 1

# ---------------------------------------------------------------------
# 2. The "last five rows" shaded-highlight formatting (style index 4 in
#    the original file) moves down as new rows are appended: rows 50-56
#    revert to the normal look, rows 57-61 become the new highlighted
#    tail. Re-apply using PasteSpecial so the existing fill/border/
#    number-format combination carries across exactly.
# ---------------------------------------------------------------------
$ws.Range("B49:F49").Copy() | Out-Null
$ws.Range("B50:F56").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Append six new daily rows (56-61)
# ---------------------------------------------------------------------
$newRows = @(
  @(56, 43945, 9, 2, 0, 0, 1),
  @(57, 43946, 2, 5, 0, 0, 2),
  @(58, 43947, 2, 2, 0, 1, 0),
  @(59, 43948, 1, 3, 0, 0, 1),
  @(60, 43949, 0, 1, 0, 2, 2),
  @(61, 43950, 0, 1, 0, 1, 0)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]

    $prev = $row - 1
    $ws.Cells.Item($row, 7).Formula  = "=B$row+G$prev"
    $ws.Cells.Item($row, 8).Formula  = "=C$row+H$prev"
    $ws.Cells.Item($row, 9).Formula  = "=D$row+I$prev"
    $ws.Cells.Item($row, 10).Formula = "=E$row+J$prev"
    $ws.Cells.Item($row, 11).Formula = "=F$row+K$prev"
    $ws.Cells.Item($row, 12).Formula = "=SUM(B$row" + ":F$row)"
    $ws.Cells.Item($row, 13).Formula = "=SUM(G$row" + ":K$row)"
    $ws.Cells.Item($row, 14).Formula = "=SUM(B$row,C$row,E$row)"
    $ws.Cells.Item($row, 15).Formula = "=SUM(G$row,H$row,J$row)"
}

# Row 56 keeps the un-shaded style (same as rows above it); rows 57-61
# get the shaded "recent" highlight that used to sit on rows 50-54.
$ws.Range("B56:F56").Copy() | Out-Null
$ws.Range("B57:F61").PasteSpecial(-4122) | Out-Null   # xlPasteFormats placeholder, overwritten below
$excel.CutCopyMode = $false

$ws.Range("B54:F54").Copy() | Out-Null
$ws.Range("B57:F61").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (shaded style)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Selection / view bookkeeping to match where the editor left off
# ---------------------------------------------------------------------
$ws.Range("S21").Select() | Out-Null
